# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a plain text value (matches the source data,
    # which stores numeric-looking prices/percentages as strings) instead of
    # letting Excel auto-convert number-like text into a floating point number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.494.26"
Set-TextValue $ws.Range("E2") "  +2.79%  "
Set-TextValue $ws.Range("D3") "2.314.88"
Set-TextValue $ws.Range("E3") "  +1.78%  "
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "310.86"
Set-TextValue $ws.Range("E5") "  +0.64%  "
Set-TextValue $ws.Range("D6") "104.40"
Set-TextValue $ws.Range("E6") "  +6.90%  "
Set-TextValue $ws.Range("E7") "  +1.28%  "
Set-TextValue $ws.Range("E8") "  +0.07%  "
Set-TextValue $ws.Range("E9") "  +8.28%  "
Set-TextValue $ws.Range("D10") "36.77"
Set-TextValue $ws.Range("E10") "  +4.94%  "
Set-TextValue $ws.Range("D11") "52.77"
Set-TextValue $ws.Range("E11") "  +1.41%  "
Set-TextValue $ws.Range("D12") "0.0813"
Set-TextValue $ws.Range("E12") "  -0.23%  "
Set-TextValue $ws.Range("E13") "  -1.15%  "
Set-TextValue $ws.Range("D14") "7.02"
Set-TextValue $ws.Range("E14") "  +2.67%  "
Set-TextValue $ws.Range("D15") "2.673.38"
Set-TextValue $ws.Range("E15") "  +1.78%  "
Set-TextValue $ws.Range("D16") "15.19"
Set-TextValue $ws.Range("E16") "  +3.91%  "
Set-TextValue $ws.Range("D17") "2.320.74"
Set-TextValue $ws.Range("E17") "  +1.92%  "
Set-TextValue $ws.Range("D18") "0.812"
Set-TextValue $ws.Range("E18") "  +2.90%  "
Set-TextValue $ws.Range("D19") "43.397.76"
Set-TextValue $ws.Range("E19") "  +2.85%  "
Set-TextValue $ws.Range("D20") "12.20"
Set-TextValue $ws.Range("E20") "  -0.77%  "
Set-TextValue $ws.Range("D21") "0.0₃0926"
Set-TextValue $ws.Range("E21") "  +2.00%  "
Set-TextValue $ws.Range("D22") "6.17"
Set-TextValue $ws.Range("E22") "  +3.36%  "
Set-TextValue $ws.Range("D23") "68.25"
Set-TextValue $ws.Range("E23") "  +0.92%  "
Set-TextValue $ws.Range("D24") "242.86"
Set-TextValue $ws.Range("E24") "  +2.55%  "
Set-TextValue $ws.Range("E25") "  +2.95%  "
Set-TextValue $ws.Range("D26") "2.62"
Set-TextValue $ws.Range("E26") "  +0.89%  "
Set-TextValue $ws.Range("D27") "1.01"
Set-TextValue $ws.Range("E27") "  +0.57%  "
Set-TextValue $ws.Range("D28") "24.86"
Set-TextValue $ws.Range("E28") "  +5.43%  "
Set-TextValue $ws.Range("D29") "2.39"
Set-TextValue $ws.Range("E29") "  +12.03%  "
Set-TextValue $ws.Range("D30") "37.30"
Set-TextValue $ws.Range("E30") "  -0.24%  "
Set-TextValue $ws.Range("E31") "  +0.61%  "
Set-TextValue $ws.Range("D32") "167.71"
Set-TextValue $ws.Range("E32") "  +2.54%  "
Set-TextValue $ws.Range("E33") "  +0.89%  "
Set-TextValue $ws.Range("E34") "  +0.00%  "
Set-TextValue $ws.Range("D35") "18.40"
Set-TextValue $ws.Range("E35") "  +3.97%  "
Set-TextValue $ws.Range("D36") "2.54"
Set-TextValue $ws.Range("E36") "  +6.77%  "
Set-TextValue $ws.Range("E37") "  +1.44%  "
Set-TextValue $ws.Range("D38") "3.07"
Set-TextValue $ws.Range("E38") "  -1.09%  "
Set-TextValue $ws.Range("D39") "1.89"
Set-TextValue $ws.Range("E39") "  +3.65%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D40") "4.49"
Set-TextValue $ws.Range("E40") "  +7.67%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.106"
Set-TextValue $ws.Range("E41") "  +1.97%  "
Set-TextValue $ws.Range("E42") "  +0.75%  "
Set-TextValue $ws.Range("D43") "2.70"
Set-TextValue $ws.Range("E43") "  +19.41%  "
Set-TextValue $ws.Range("D44") "0.0294"
Set-TextValue $ws.Range("E44") "  +3.70%  "
Set-TextValue $ws.Range("D45") "1.993.85"
Set-TextValue $ws.Range("E45") "  +2.49%  "
Set-TextValue $ws.Range("E46") "  +1.82%  "
Set-TextValue $ws.Range("D47") "3.07"
Set-TextValue $ws.Range("E47") "  +3.72%  "
Set-TextValue $ws.Range("E48") "  +2.35%  "
Set-TextValue $ws.Range("D49") "57.15"
Set-TextValue $ws.Range("E49") "  +5.73%  "
Set-TextValue $ws.Range("E50") "  +1.76%  "
Set-TextValue $ws.Range("E51") "  +9.32%  "
